$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add a new journal entry on row 6 (Gabriel Pereira - Creation de maquette)
# ------------------------------------------------------------------

# 1) Start from row 4's formatting (plain "middle of table" style already
#    used by row 6 for columns C:E and H:I) then overlay row 5's styling
#    for the columns (A,B,F,G) that sit right below row 5's bottom border
#    so no redundant top border is drawn (mirrors how Excel behaves when
#    you fill an adjoining row down from the row above).
$ws.Range("A4:I4").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A5:B5").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F5:G5").Copy()
$ws.Range("F6:G6").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# 2) Fill in the new row's data
$ws.Range("A6").Value = "Gabriel Pereira"
$ws.Range("B6").Value = 43896
$ws.Range("C6").Value = 0.95347222222222217
$ws.Range("D6").Value = 0.99444444444444446
$ws.Range("E6").Formula = "=D6-C6"
$ws.Range("F6").Value = "draw.io"
$ws.Range("G6").Value = "Bataille Navale"
$ws.Range("H6").Value = "Création de maquette"
$ws.Range("I6").Value = "Création d'une maquette pour la bataille navale"

# 3) Match the author's final selection in the sheet
$ws.Range("D7").Select()
